$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","E","F","G","H","I","J","K")

$data = @(
    @(9.864663296621467, 5.125298282309029, 21.80972674123574, 42.38265513367197, 34.01192259201394, 15.65872118919969, 22.49605537051562, 8.418107799487522, 9.993832589554678),
    @(9.55549872407156, 4.895355888373637, 21.48453272885329, 42.21212241171112, 34.1556105386508, 15.72664154162027, 22.61967449578849, 8.44615483969379, 9.78612181453529),
    @(9.362095983091672, 4.74940854058992, 21.28722839087556, 42.11952840203458, 34.25653226181844, 15.77137426321582, 22.70064455172188, 8.464416462934739, 9.658065727796902),
    @(9.282508170836159, 4.688823276460693, 21.20751764228121, 42.0848667056283, 34.30082806960782, 15.7903640744783, 22.73491293680597, 8.472120305955128, 9.605823566795648),
    @(9.26924970017501, 4.678698838765654, 21.19432637912421, 42.07929728799849, 34.30837408170555, 15.79356324697488, 22.74067997627473, 8.473415364553752, 9.597147438650376),
    @(9.361025599556486, 4.74859583952905, 21.28615045245737, 42.11904847785767, 34.25711684932662, 15.77162728693278, 22.70110155751247, 8.464519297946214, 9.657361310201393),
    @(9.758876854190438, 5.04705095685383, 21.6971674068017, 42.32135758516802, 34.0588171997228, 15.68151084513216, 22.53762639775481, 8.427562682560277, 9.922364994381248),
    @(10.50572262353781, 5.591294906350091, 22.51764625786159, 42.81289806644135, 33.77170025134124, 15.5288695473368, 22.25734226281344, 8.363329202214734, 10.43482334724291),
    @(11.02811766129935, 5.962513985421724, 23.12355322785306, 43.22979972295032, 33.62409916053475, 15.43145760946246, 22.07609499010288, 8.321133313667021, 10.80299858585328),
    @(11.25897955839598, 6.124585578724298, 23.39868352077591, 43.43105059859537, 33.57096676691162, 15.39035433517651, 21.99902625060979, 8.303016889368022, 10.96791539465878),
    @(11.34534885075231, 6.184941355694837, 23.50270050611093, 43.50887783302422, 33.55288147664191, 15.37525235540532, 21.97061874076153, 8.296311409368103, 11.02993677683523),
    @(11.32679567283768, 6.171988464686001, 23.48030760889061, 43.49204522398307, 33.55668566666035, 15.37848422566904, 21.97670221295616, 8.297748673591929, 11.01659935430348),
    @(11.26610673642262, 6.129571676045529, 23.40724498976573, 43.43742136640515, 33.56943798961446, 15.3891026003454, 21.99667356352695, 8.302462124408944, 10.97302685814939),
    @(11.2287936853949, 6.103456627935605, 23.36246720718361, 43.40417182696645, 33.5775147139172, 15.39566699051984, 22.00900783571273, 8.30536940215214, 10.94627986547758),
    @(11.01288699977784, 5.951781893256398, 23.105554394486, 43.21687671724663, 33.62785517405504, 15.4342085118754, 22.08124014844416, 8.322338944979119, 10.79216401381589),
    @(10.8786406424861, 5.856962920882157, 22.94774869351985, 43.10491370127401, 33.6623407848267, 15.45867557575582, 22.12693269751296, 8.333025274295872, 10.69691747293828),
    @(10.80079181089035, 5.801787638106967, 22.85694160512945, 43.04161063427284, 33.68349376780878, 15.473050431743, 22.15372020068537, 8.339273317535016, 10.64189580688287),
    @(10.77432719548704, 5.782997989783683, 22.82619177731008, 43.02036686142348, 33.69088150693479, 15.47796934724814, 22.16287687025673, 8.341406248077522, 10.62322735916968),
    @(10.89299759549159, 5.867122892588125, 22.96455237958605, 43.11671936967672, 33.65853320017368, 15.45603974444443, 22.12201622094456, 8.331877188521615, 10.70708174588478),
    @(11.28396171141454, 6.142058406746165, 23.42871059282731, 43.45342222489722, 33.56563694747314, 15.38597115130961, 21.99078638942436, 8.301073470295659, 10.98583725745088),
    @(11.53331064805924, 6.315801462522892, 23.73103777634266, 43.68288176216615, 33.51679251116764, 15.34287628712936, 21.90954901721942, 8.281843718809013, 11.16549438163063),
    @(11.40081658476417, 6.223626825972441, 23.56980524458067, 43.55957178573868, 33.54176949439561, 15.36562941003632, 21.95249155283706, 8.292024538279273, 11.06985783798088),
    @(10.88650889418507, 5.862531630135772, 22.95695568262493, 43.11137870267562, 33.66025047783775, 15.45723044313041, 22.12423734582003, 8.3323959130901, 10.70248729316637),
    @(10.30790581712485, 5.448865366755222, 22.29473066079392, 42.66996646355022, 33.83834052131955, 15.56757984691974, 22.32884167094844, 8.37982680068351, 10.29737998514642)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = "$($cols[$j])$rowNum"
        $ws.Range($addr).Value = $rowVals[$j]
    }
}

Write-Output "Updated $($data.Length) rows"
